# Generate Report for Archive
# - Update status text "Ready for handoff" -> "In Translation" everywhere it appears
# - Narrow the status-related columns (width 17.2159881591797 -> 13.4101845877511)

$wb = $excel.ActiveWorkbook

# --- Update status cell values on every sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("E4").Value = "In Translation"
$overview.Range("F4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"

# --- Resize the Status-related columns to their new (narrower) width ---
# Target stored width is 13.4101845877511 characters; the ColumnWidth value
# below is the input that this runtime's column-width quantization maps
# closest to that target.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
